$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.5941044984698646
$ws.Range("E2").Value = 0.5941044984698646

# Row 3
$ws.Range("D3").Value = 0.001144904478007686
$ws.Range("E3").Value = 0.001144904478007686

# Row 4
$ws.Range("D4").Value = 0.9488775416732173
$ws.Range("E4").Value = 0.9488775416732173

# Row 5
$ws.Range("D5").Value = 0.04907159283619292
$ws.Range("E5").Value = 0.04907159283619292

# Row 6
$ws.Range("D6").Value = 0.2303919460329533
$ws.Range("E6").Value = 0.2303919460329533

# Row 7
$ws.Range("D7").Value = 0.8534811682151161
$ws.Range("E7").Value = 0.1465188317848839

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.001207019421165155
$ws.Range("E8").Value = 0.9987929805788348

# Row 9
$ws.Range("D9").Value = 0.9868877157573662
$ws.Range("E9").Value = 0.01311228424263378

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 0.0007363213694898383
$ws.Range("E10").Value = 0.9992636786305101

# Row 11
$ws.Range("D11").Value = 0.9999999381780804
$ws.Range("E11").Value = [double]"6.18219195835934E-08"
$ws.Range("F11").Value = 1.829360246658325
$ws.Range("G11").Value = 0.6

# Row 12
$ws.Range("D12").Value = 0.8562582463161288
$ws.Range("E12").Value = 0.8562582463161288

# Row 13
$ws.Range("D13").Value = [double]"3.684100776501542E-05"
$ws.Range("E13").Value = [double]"3.684100776501542E-05"

# Row 14
$ws.Range("D14").Value = 0.9974794082974022
$ws.Range("E14").Value = 0.9974794082974022

# Row 15
$ws.Range("D15").Value = 0.004237297386042994
$ws.Range("E15").Value = 0.004237297386042994

# Row 16
$ws.Range("D16").Value = 0.208829455497294
$ws.Range("E16").Value = 0.208829455497294

# Row 17
$ws.Range("D17").Value = 0.9882871501413207
$ws.Range("E17").Value = 0.01171284985867926

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = [double]"2.532348319467842E-09"
$ws.Range("E18").Value = 0.9999999974676517

# Row 19
$ws.Range("D19").Value = 0.9570022233664972
$ws.Range("E19").Value = 0.04299777663350279

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = [double]"2.12927158470294E-05"
$ws.Range("E20").Value = 0.9999787072841529

# Row 21
$ws.Range("D21").Value = 0.9999999999999976
$ws.Range("E21").Value = [double]"2.442490654175344E-15"
$ws.Range("F21").Value = 3.876851558685303
$ws.Range("G21").Value = 0.6
